$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(12, 8).Value = 920.7778
$ws.Cells.Item(12, 9).Value = 633.1667
$ws.Cells.Item(12, 10).Value = 1496
$ws.Cells.Item(12, 11).Value = 633.1667
$ws.Cells.Item(12, 12).Value = 1496
$ws.Cells.Item(12, 13).Value = -463.1667
$ws.Cells.Item(12, 14).Value = -1836

$ws.Cells.Item(18, 8).Value = 15922.75
$ws.Cells.Item(18, 9).Value = 19716.666
$ws.Cells.Item(18, 11).Value = 19716.666
$ws.Cells.Item(18, 13).Value = -19432.666

$ws.Cells.Item(19, 8).Value = 449.33334
$ws.Cells.Item(19, 9).Value = 333
$ws.Cells.Item(19, 11).Value = 333
$ws.Cells.Item(19, 13).Value = -158

$ws.Cells.Item(33, 8).Value = 318.4
$ws.Cells.Item(33, 9).Value = 323
$ws.Cells.Item(33, 11).Value = 323
$ws.Cells.Item(33, 13).Value = -94

$ws.Cells.Item(55, 8).Value = 641.4375
$ws.Cells.Item(55, 9).Value = 405.2857
$ws.Cells.Item(55, 10).Value = 825.1111
$ws.Cells.Item(55, 11).Value = 405.2857
$ws.Cells.Item(55, 12).Value = 825.1111
$ws.Cells.Item(55, 13).Value = -191.2857
$ws.Cells.Item(55, 14).Value = -1253.1111

$ws.Cells.Item(70, 8).Value = 3826.647
$ws.Cells.Item(70, 10).Value = 3898.4546
$ws.Cells.Item(70, 12).Value = 11695.3638
$ws.Cells.Item(70, 14).Value = -12235.3638

$ws.Cells.Item(73, 8).Value = 3826.647
$ws.Cells.Item(73, 10).Value = 3898.4546
$ws.Cells.Item(73, 12).Value = 11695.3638
$ws.Cells.Item(73, 14).Value = -13567.3638

$ws.Cells.Item(74, 8).Value = 3887.7058
$ws.Cells.Item(74, 9).Value = 2917.3635
$ws.Cells.Item(74, 11).Value = 2917.3635
$ws.Cells.Item(74, 13).Value = -1981.3635

$ws.Cells.Item(77, 8).Value = 3887.7058
$ws.Cells.Item(77, 9).Value = 2917.3635
$ws.Cells.Item(77, 11).Value = 14586.8175
$ws.Cells.Item(77, 13).Value = -9906.817499999999

$ws.Cells.Item(135, 8).Value = 1532.25
$ws.Cells.Item(135, 9).Value = 1706.95
$ws.Cells.Item(135, 11).Value = 15362.55
$ws.Cells.Item(135, 13).Value = -12827.55

$ws.Cells.Item(137, 8).Value = 1510.8334
$ws.Cells.Item(137, 9).Value = 1366.25
$ws.Cells.Item(137, 11).Value = 4098.75
$ws.Cells.Item(137, 13).Value = -1548.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 992.1556
$ws.Cells.Item(2, 9).Value = 959.7632
$ws.Cells.Item(2, 11).Value = 959.7632
$ws.Cells.Item(2, 13).Value = -846.7632

$ws.Cells.Item(32, 8).Value = 3781.5264
$ws.Cells.Item(32, 9).Value = 3195.3518
$ws.Cells.Item(32, 11).Value = 3195.3518
$ws.Cells.Item(32, 13).Value = -2908.3518

$ws.Cells.Item(37, 8).Value = 4567231
$ws.Cells.Item(37, 9).Value = 25007316
$ws.Cells.Item(37, 10).Value = 24989.445
$ws.Cells.Item(37, 11).Value = 25007316
$ws.Cells.Item(37, 12).Value = 24989.445
$ws.Cells.Item(37, 13).Value = -25007043
$ws.Cells.Item(37, 14).Value = -25535.445

$ws.Cells.Item(45, 8).Value = 6964.6
$ws.Cells.Item(45, 9).Value = 9656.691999999999
$ws.Cells.Item(45, 11).Value = 9656.691999999999
$ws.Cells.Item(45, 13).Value = -9279.691999999999

$ws.Cells.Item(88, 8).Value = 100000150
$ws.Cells.Item(88, 9).Value = 300
$ws.Cells.Item(88, 11).Value = 300
$ws.Cells.Item(88, 13).Value = 106

$ws.Cells.Item(91, 8).Value = 100000150
$ws.Cells.Item(91, 9).Value = 300
$ws.Cells.Item(91, 11).Value = 300
$ws.Cells.Item(91, 13).Value = 1104

$ws.Cells.Item(111, 8).Value = 53977
$ws.Cells.Item(111, 10).Value = 53977
$ws.Cells.Item(111, 12).Value = 53977
$ws.Cells.Item(111, 14).Value = -62157

$ws.Cells.Item(116, 8).Value = 992.1556
$ws.Cells.Item(116, 9).Value = 959.7632
$ws.Cells.Item(116, 11).Value = 959.7632
$ws.Cells.Item(116, 13).Value = 1334.2368

$ws.Cells.Item(122, 8).Value = 2105.9119
$ws.Cells.Item(122, 9).Value = 1375.8
$ws.Cells.Item(122, 10).Value = 3148.9285
$ws.Cells.Item(122, 11).Value = 4127.4
$ws.Cells.Item(122, 12).Value = 9446.7855
$ws.Cells.Item(122, 13).Value = -1677.4
$ws.Cells.Item(122, 14).Value = -14346.7855

$ws.Cells.Item(132, 8).Value = 2213.2727
$ws.Cells.Item(132, 9).Value = 2189.1316
$ws.Cells.Item(132, 10).Value = 2366.1667
$ws.Cells.Item(132, 11).Value = 6567.3948
$ws.Cells.Item(132, 12).Value = 7098.500100000001
$ws.Cells.Item(132, 13).Value = -4037.3948
$ws.Cells.Item(132, 14).Value = -12158.5001

$ws.Cells.Item(140, 8).Value = 120342.836
$ws.Cells.Item(140, 10).Value = 120342.836
$ws.Cells.Item(140, 12).Value = 120342.836
$ws.Cells.Item(140, 14).Value = -130702.836

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 992.1556
$ws.Cells.Item(3, 9).Value = 959.7632
$ws.Cells.Item(3, 11).Value = 959.7632
$ws.Cells.Item(3, 13).Value = -845.7632

$ws.Cells.Item(64, 8).Value = 987.9
$ws.Cells.Item(64, 9).Value = 777.5
$ws.Cells.Item(64, 10).Value = 1303.5
$ws.Cells.Item(64, 11).Value = 777.5
$ws.Cells.Item(64, 12).Value = 1303.5
$ws.Cells.Item(64, 13).Value = -552.5
$ws.Cells.Item(64, 14).Value = -1753.5

$ws.Cells.Item(67, 8).Value = 987.9
$ws.Cells.Item(67, 9).Value = 777.5
$ws.Cells.Item(67, 10).Value = 1303.5
$ws.Cells.Item(67, 11).Value = 777.5
$ws.Cells.Item(67, 12).Value = 1303.5
$ws.Cells.Item(67, 13).Value = 2.5
$ws.Cells.Item(67, 14).Value = -2863.5

$ws.Cells.Item(74, 8).Value = 0
$ws.Cells.Item(74, 10).Value = 0
$ws.Cells.Item(74, 12).Value = 0
$ws.Cells.Item(74, 14).ClearContents()

$ws.Cells.Item(77, 8).Value = 0
$ws.Cells.Item(77, 10).Value = 0
$ws.Cells.Item(77, 12).Value = 0
$ws.Cells.Item(77, 14).ClearContents()

$ws.Cells.Item(99, 8).Value = 24595.334
$ws.Cells.Item(99, 9).Value = 27507.25
$ws.Cells.Item(99, 11).Value = 27507.25
$ws.Cells.Item(99, 13).Value = -26009.25

$ws.Cells.Item(134, 8).Value = 1540.2632
$ws.Cells.Item(134, 9).Value = 1540.2632
$ws.Cells.Item(134, 11).Value = 4620.7896
$ws.Cells.Item(134, 13).Value = -2085.7896

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 19230.684
$ws.Cells.Item(31, 9).Value = 3027.818
$ws.Cells.Item(31, 10).Value = 41509.625
$ws.Cells.Item(31, 11).Value = 3027.818
$ws.Cells.Item(31, 12).Value = 41509.625
$ws.Cells.Item(31, 13).Value = -2732.818
$ws.Cells.Item(31, 14).Value = -42099.625

$ws.Cells.Item(34, 8).Value = 19230.684
$ws.Cells.Item(34, 9).Value = 3027.818
$ws.Cells.Item(34, 10).Value = 41509.625
$ws.Cells.Item(34, 11).Value = 3027.818
$ws.Cells.Item(34, 12).Value = 41509.625
$ws.Cells.Item(34, 13).Value = -2825.818
$ws.Cells.Item(34, 14).Value = -41913.625

$ws.Cells.Item(132, 8).Value = 3742.2273
$ws.Cells.Item(132, 9).Value = 4141.3887
$ws.Cells.Item(132, 11).Value = 12424.1661
$ws.Cells.Item(132, 13).Value = -9894.166100000002

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(98, 8).Value = 0
$ws.Cells.Item(98, 10).Value = 0
$ws.Cells.Item(98, 12).Value = 0
$ws.Cells.Item(98, 14).ClearContents()

$ws.Cells.Item(121, 8).Value = 20897554
$ws.Cells.Item(121, 10).Value = 114069.445
$ws.Cells.Item(121, 12).Value = 342208.335
$ws.Cells.Item(121, 14).Value = -344828.335

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(74, 8).Value = 48131
$ws.Cells.Item(74, 10).Value = 48131
$ws.Cells.Item(74, 12).Value = 48131
$ws.Cells.Item(74, 14).Value = -50003

$ws.Cells.Item(77, 8).Value = 48131
$ws.Cells.Item(77, 10).Value = 48131
$ws.Cells.Item(77, 12).Value = 144393
$ws.Cells.Item(77, 14).Value = -153753

$ws.Cells.Item(80, 8).Value = 3063.8572
$ws.Cells.Item(80, 9).Value = 2492
$ws.Cells.Item(80, 10).Value = 6495
$ws.Cells.Item(80, 11).Value = 2492
$ws.Cells.Item(80, 12).Value = 6495
$ws.Cells.Item(80, 13).Value = -1494
$ws.Cells.Item(80, 14).Value = -8491

$ws.Cells.Item(83, 8).Value = 3063.8572
$ws.Cells.Item(83, 9).Value = 2492
$ws.Cells.Item(83, 10).Value = 6495
$ws.Cells.Item(83, 11).Value = 12460
$ws.Cells.Item(83, 12).Value = 32475
$ws.Cells.Item(83, 13).Value = -7468
$ws.Cells.Item(83, 14).Value = -42459

$ws.Cells.Item(122, 8).Value = 2900.7576
$ws.Cells.Item(122, 9).Value = 2819.2856
$ws.Cells.Item(122, 10).Value = 3357
$ws.Cells.Item(122, 11).Value = 8457.856800000001
$ws.Cells.Item(122, 12).Value = 10071
$ws.Cells.Item(122, 13).Value = -6007.856800000001
$ws.Cells.Item(122, 14).Value = -14971

$ws.Cells.Item(132, 8).Value = 2356.6875
$ws.Cells.Item(132, 9).Value = 2406.9312
$ws.Cells.Item(132, 10).Value = 1871
$ws.Cells.Item(132, 11).Value = 7220.7936
$ws.Cells.Item(132, 12).Value = 5613
$ws.Cells.Item(132, 13).Value = -4690.7936
$ws.Cells.Item(132, 14).Value = -10673

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 54827
$ws.Cells.Item(46, 9).Value = 212561.5
$ws.Cells.Item(46, 10).Value = 2248.8333
$ws.Cells.Item(46, 11).Value = 212561.5
$ws.Cells.Item(46, 12).Value = 2248.8333
$ws.Cells.Item(46, 13).Value = -212373.5
$ws.Cells.Item(46, 14).Value = -2624.8333

$ws.Cells.Item(115, 8).Value = 0
$ws.Cells.Item(115, 10).Value = 0
$ws.Cells.Item(115, 12).Value = 0
$ws.Cells.Item(115, 14).ClearContents()

$ws.Cells.Item(118, 8).Value = 47000
$ws.Cells.Item(118, 10).Value = 47000
$ws.Cells.Item(118, 12).Value = 47000
$ws.Cells.Item(118, 14).Value = -50314

$ws.Cells.Item(132, 8).Value = 3501.5278
$ws.Cells.Item(132, 9).Value = 3365.7083
$ws.Cells.Item(132, 10).Value = 3773.1667
$ws.Cells.Item(132, 11).Value = 10097.1249
$ws.Cells.Item(132, 12).Value = 11319.5001
$ws.Cells.Item(132, 13).Value = -7567.124899999999
$ws.Cells.Item(132, 14).Value = -16379.5001

$ws.Cells.Item(136, 8).Value = 3829.7
$ws.Cells.Item(136, 9).Value = 1976.75
$ws.Cells.Item(136, 10).Value = 5065
$ws.Cells.Item(136, 11).Value = 5930.25
$ws.Cells.Item(136, 12).Value = 15195
$ws.Cells.Item(136, 13).Value = -3380.25
$ws.Cells.Item(136, 14).Value = -20295

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(126, 8).Value = 1558.04
$ws.Cells.Item(126, 10).Value = 1742.1428
$ws.Cells.Item(126, 12).Value = 5226.428400000001
$ws.Cells.Item(126, 14).Value = -10166.4284

$ws.Cells.Item(135, 8).Value = 91617.336
$ws.Cells.Item(135, 10).Value = 91617.336
$ws.Cells.Item(135, 12).Value = 91617.336
$ws.Cells.Item(135, 14).Value = -101757.336

$ws.Cells.Item(136, 8).Value = 2628.6428
$ws.Cells.Item(136, 9).Value = 2150.0833
$ws.Cells.Item(136, 10).Value = 5500
$ws.Cells.Item(136, 11).Value = 6450.249899999999
$ws.Cells.Item(136, 12).Value = 16500
$ws.Cells.Item(136, 13).Value = -3900.249899999999
$ws.Cells.Item(136, 14).Value = -21600
